$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 (JAMA article) - reset to "unknown"/"not found" placeholder values
$ws.Range("C2").Value = "Unknown Title"
$ws.Range("E2").Value = "[]"
$ws.Range("F2").Value = "not found"
$ws.Range("G2").Value = "N/A"
$ws.Range("H2").Value = "1970-01-01"
$ws.Range("J2").Value = ""

# Row 3 (JAMA article) - reset to "unknown"/"not found" placeholder values
$ws.Range("C3").Value = "Unknown Title"
$ws.Range("E3").Value = "[]"
$ws.Range("F3").Value = "not found"
$ws.Range("G3").Value = "N/A"
$ws.Range("H3").Value = "1970-01-01"
$ws.Range("J3").Value = ""

# Row 4 - full text authors and publisher info
$ws.Range("E4").Value = "[Ruchong%Chen%NULL%0, Wenhua%Liang%NULL%0, Mei%Jiang%NULL%0, Weijie%Guan%NULL%0, Chen%Zhan%NULL%0, Tao%Wang%NULL%0, Chunli%Tang%NULL%0, Ling%Sang%NULL%0, Jiaxing%Liu%NULL%0, Zhengyi%Ni%NULL%0, Yu%Hu%NULL%0, Lei%Liu%NULL%0, Hong%Shan%NULL%0, Chunliang%Lei%NULL%0, Yixiang%Peng%NULL%0, Li%Wei%NULL%0, Yong%Liu%NULL%0, Yahua%Hu%NULL%0, Peng%Peng%NULL%0, Jianming%Wang%NULL%0, Jiyang%Liu%NULL%0, Zhong%Chen%NULL%0, Gang%Li%NULL%0, Zhijian%Zheng%NULL%0, Shaoqin%Qiu%NULL%0, Jie%Luo%NULL%0, Changjiang%Ye%NULL%0, Shaoyong%Zhu%NULL%0, Xiaoqing%Liu%NULL%0, Linling%Cheng%NULL%0, Feng%Ye%NULL%0, Jinping%Zheng%NULL%0, Nuofu%Zhang%NULL%0, Yimin%Li%NULL%0, Jianxing%He%NULL%0, Shiyue%Li%lishiyue@188.com%0, Nanshan%Zhong%NULL%0, NULL%NULL%NULL%0]"
$ws.Range("I4").Value = ""
$ws.Range("J4").Value = "American College of Chest Physicians. Published by Elsevier Inc."

# Row 5 - full text authors and publisher info
$ws.Range("E5").Value = "[Wei-jie%Guan%NULL%0, Zheng-yi%Ni%NULL%0, Zheng-yi%Ni%NULL%0, Yu%Hu%NULL%0, Wen-hua%Liang%NULL%0, Chun-quan%Ou%NULL%0, Jian-xing%He%NULL%0, Lei%Liu%NULL%0, Hong%Shan%NULL%0, Chun-liang%Lei%NULL%0, David S.C.%Hui%NULL%0, Bin%Du%NULL%0, Lan-juan%Li%NULL%0, Guang%Zeng%NULL%0, Kwok-Yung%Yuen%NULL%0, Ru-chong%Chen%NULL%0, Chun-li%Tang%NULL%0, Tao%Wang%NULL%0, Ping-yan%Chen%NULL%0, Jie%Xiang%NULL%0, Shi-yue%Li%NULL%0, Jin-lin%Wang%NULL%0, Zi-jing%Liang%NULL%0, Yi-xiang%Peng%NULL%0, Li%Wei%NULL%0, Yong%Liu%NULL%0, Ya-hua%Hu%NULL%0, Peng%Peng%NULL%0, Jian-ming%Wang%NULL%0, Ji-yang%Liu%NULL%0, Zhong%Chen%NULL%0, Gang%Li%NULL%0, Zhi-jian%Zheng%NULL%0, Shao-qin%Qiu%NULL%0, Jie%Luo%NULL%0, Chang-jiang%Ye%NULL%0, Shao-yong%Zhu%NULL%0, Nan-shan%Zhong%NULL%0]"
$ws.Range("I5").Value = ""
$ws.Range("J5").Value = "Massachusetts Medical Society"
